$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.284.13"
$ws.Range("E2").Value = "  +1.92%  "

$ws.Range("D3").Value = "2.061.73"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'233.01"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  +2.63%  "

$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.86"
$ws.Range("E7").Value = "  +4.41%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'0.383"
$ws.Range("E9").Value = "  +2.71%  "

$ws.Range("D10").Value = "'58.02"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("D11").Value = "'0.0759"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("D13").Value = "2.359.77"
$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").Value = "'14.37"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "'20.74"
$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("D16").Value = "'0.774"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "2.054.73"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("D19").Value = "37.193.32"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").Value = "'6.34"
$ws.Range("E20").Value = "  +13.36%  "

$ws.Range("D21").Value = "'69.11"
$ws.Range("E21").Value = "  +1.86%  "

$ws.Range("D22").Value = "0.0₃0811"
$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("D23").Value = "'225.06"
$ws.Range("E23").Value = "  +1.61%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  +2.16%  "

$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").Value = "'165.80"
$ws.Range("E27").Value = "  +1.78%  "

$ws.Range("E28").Value = "  +7.42%  "

$ws.Range("D29").Value = "'8.80"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("E30").Value = "  -2.93%  "

$ws.Range("D31").Value = "'19.08"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("D33").Value = "'4.47"
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").Value = "'0.0614"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("D35").Value = "'2.54"
$ws.Range("E35").Value = "  +3.92%  "

$ws.Range("D36").Value = "'4.53"
$ws.Range("E36").Value = "  +6.13%  "

$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'5.83"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("E39").Value = "  -0.92%  "

$ws.Range("E40").Value = "  -1.48%  "

$ws.Range("D41").Value = "'4.58"
$ws.Range("E41").Value = "  +13.45%  "

$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").Value = "1.488.15"
$ws.Range("E43").Value = "  +1.42%  "

$ws.Range("D44").Value = "'96.72"
$ws.Range("E44").Value = "  +2.24%  "

$ws.Range("E45").Value = "  +3.96%  "

$ws.Range("E46").Value = "  -1.84%  "

$ws.Range("E47").Value = "  +3.06%  "

$ws.Range("D48").Value = "'15.35"
$ws.Range("E48").Value = "  -2.07%  "

$ws.Range("E49").Value = "  +1.64%  "

$ws.Range("E50").Value = "  +4.14%  "

$ws.Range("E51").Value = "  +2.12%  "
